$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (QBTS) - numeric recalculation updates
$ws.Range("D2").Value = 27.09
$ws.Range("E2").Value = 59.5
$ws.Range("F2").Value = 19.5
$ws.Range("N2").Value = 50.68470204858703

# Row 3 becomes IonQ, Inc. / IONQ (was International Business Machines / IBM)
$ws.Range("B3").Value = "IonQ, Inc."
$ws.Range("C3").Value = "IONQ"
$ws.Range("D3").Value = 52.01
$ws.Range("E3").Value = 57.6
$ws.Range("F3").Value = 5.5
$ws.Range("H3").Value = 53
$ws.Range("I3").Value = 60
$ws.Range("J3").Value = 70
$ws.Range("K3").Value = 57.2
$ws.Range("N3").Value = 50.68470204858703

# Row 4 becomes Rigetti Computing, Inc. / RGTI (was IonQ, Inc. / IONQ)
$ws.Range("B4").Value = "Rigetti Computing, Inc."
$ws.Range("C4").Value = "RGTI"
$ws.Range("D4").Value = 28.06
$ws.Range("E4").Value = 56.4
$ws.Range("F4").Value = 9.72
$ws.Range("G4").Value = 50
$ws.Range("H4").Value = 60
$ws.Range("I4").Value = 66
$ws.Range("J4").Value = 83
$ws.Range("K4").Value = 56.6
$ws.Range("N4").Value = 50.68470204858703

# Row 5 becomes International Business Machines / IBM (was Rigetti Computing, Inc. / RGTI)
$ws.Range("B5").Value = "International Business Machines"
$ws.Range("C5").Value = "IBM"
$ws.Range("D5").Value = 308.48
$ws.Range("E5").Value = 52.6
$ws.Range("F5").Value = -0.03
$ws.Range("G5").Value = 40
$ws.Range("J5").Value = 63
$ws.Range("K5").Value = 53.6
$ws.Range("N5").Value = 50.68470204858703
